$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("Dahdahle", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "Dahdaleh", 2)
        }
    }
}
